$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games), rows 2-13, columns B-E and G.
# F (Win) column is left unchanged.
$data = @{
    2  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569;  G = 6.048734245549538 }
    3  = @{ B = 1.505614041169197;  C = 0.3375848360084654; D = 3.082599426703578;   E = 6.48142807727062;    G = 11.40722638115186 }
    4  = @{ B = 0.7287194209349384; C = 1.65323645889881;   D = 0.1529057820181812;  E = 0.4998867070740569;  G = 3.034748368925986 }
    5  = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569;  G = 4.371470058157054 }
    6  = @{ B = 0.06328177979961902;C = 0.05231270169004087;D = 0.1529057820181812;  E = 6.48142807727062;    G = 6.749928340778461 }
    7  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 16.98373111632243;   E = 0.4998867070740569;  G = 22.31973251085698 }
    8  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.1529057820181812;  E = 0.4998867070740569;  G = 5.488907176552729 }
    9  = @{ B = 0.1554434735375247; C = 0.3375848360084654; D = 3.082599426703578;   E = 0.4998867070740569;  G = 4.075514443323626 }
    10 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569;  G = 6.048734245549538 }
    11 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 0.7127328510149897;  E = 6.48142807727062;    G = 10.35301142835362 }
    12 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 157.8057217802531;   E = 6.48142807727062;    G = 167.4460003575917 }
    13 = @{ B = 0.7287194209349384; C = 1.65323645889881;   D = 0.1529057820181812;  E = 0.4998867070740569;  G = 3.034748368925986 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
